$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new job-tracker row (row 65) with the Ericsson ASIC Architect posting.
$ws.Cells.Item(65, 1).Value = "Ericsson"
$ws.Cells.Item(65, 2).Value = "Sweden"
$ws.Cells.Item(65, 3).Value = "Lund"
$ws.Cells.Item(65, 4).Value = "ASIC Architect"
$ws.Cells.Item(65, 5).Value = 768407
$ws.Cells.Item(65, 6).Value = "https://jobs.ericsson.com/careers/job/563121765371416"
$ws.Rows.Item(65).RowHeight = 15.9

# Turn the URL cell into a real hyperlink, matching the style used by the rest of column F.
$ws.Hyperlinks.Add($ws.Cells.Item(65, 6), "https://jobs.ericsson.com/careers/job/563121765371416")
$ws.Cells.Item(65, 6).Style = "Hyperlink"

# Move the selection to reflect where the author ended up after adding the row.
$ws.Range("E66").Select()
